$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cheque field content updates ---

# Date (H2)
$ws.Range("H2").Value = "SEPTEMBER 01, 2019"

# Payee name (B3)
$ws.Range("B3").Value = "EMPERADOR DISTILLERS, INC."

# Amount (H3) - the cell already carries a Text ("@") number format, so this
# is stored as text, matching the new shared-string-backed value.
$ws.Range("H3").Value = "227,900.00"

# Amount in words (B4)
$ws.Range("B4").Value = "TWO HUNDRED TWENTY-SEVEN THOUSAND NINE HUNDRED ONLY"

# --- Formatting updates ---

# The amount cell is now center-aligned instead of right-aligned.
$ws.Range("H3").HorizontalAlignment = -4108  # xlCenter

# Row heights shrank slightly for rows 2-4.
$ws.Rows.Item(2).RowHeight = 19.2
$ws.Rows.Item(3).RowHeight = 19.8
$ws.Rows.Item(4).RowHeight = 19.8

# Column widths: column A widened a touch, and column G split off from the
# B:G block into its own (slightly wider) width.
$ws.Columns.Item(1).ColumnWidth = 9.15
$ws.Columns.Item(7).ColumnWidth = 10.3

# --- Selection state ---
$ws.Range("H3").Select()
